# "made final tweaks David wanted"
#
# - Row 10 used to hold a stray time value in C10 (0.0764 ~= 01:50 AM).
#   That was replaced with a proper timesheet entry: date in A10, hours
#   in B10, and C10 cleared out (but keeps its time-formatted style).
# - A brand new row 11 was added with the next timesheet entry
#   (date in A11, hours in B11, C11 left blank with the time style).
# - B18's SUM(B2:B15) total recalculates automatically to include the
#   two new entries.
# - Columns were tidied up (auto-fit) and the selection left on the
#   total cell, B18.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10: turn the old stray C10 time value into a real date/hours entry.
$ws.Range("C10").ClearContents()
$ws.Range("A10").Value = 42937
$ws.Range("A10").NumberFormat = "MM/DD/YY"
$ws.Range("B10").Value = 1

# Row 11: new timesheet entry added underneath.
$ws.Range("A11").Value = 42940
$ws.Range("A11").NumberFormat = "MM/DD/YY"
$ws.Range("B11").Value = 0.1
$ws.Range("C11").NumberFormat = "HH:MM:SS\ AM/PM"

# Tidy up column widths now that the data has changed.
$ws.Columns.AutoFit()

# Leave the selection on the (recalculated) total cell.
$ws.Range("B18").Select()
